$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Rename header G1 "play" -> "play_vslh" (vs left-handed pitching) and add a new
# header H1 "play_vsrh" (vs right-handed pitching) so the two splits get their own column.
$ws.Range("G1").Value = "play_vslh"
$ws.Range("H1").Value = "play_vsrh"

# Seed the new vsrh column with the same play outcome already recorded in the
# vslh column for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 7).Value2
}

# Update view/selection to match post-edit state (scrolled one column over, new
# column highlighted)
$ws.Range("H2:H" + $lastRow).Select() | Out-Null
